# Update model_comparison sheet: add "PR AUC" metric rows for each model
# and resort the data (Accuracy, Precision, Recall, ROC AUC, PR AUC) per model.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data table (excluding header row). Columns: Model, Metric, Value
$data = @(
    @("RBD-HA",     "Accuracy",  0.8694714308),
    @("RBD-HA",     "Precision", 0.89020770790000003),
    @("RBD-HA",     "Recall",    0.92735701800000003),
    @("RBD-HA",     "ROC AUC",   0.91839247940000002),
    @("RBD-HA",     "PR AUC",    0.96001636980000005),
    @("Spike-HA",   "Accuracy",  0.82224768400000003),
    @("Spike-HA",   "Precision", 0.85643565649999998),
    @("Spike-HA",   "Recall",    0.88415670390000001),
    @("Spike-HA",   "ROC AUC",   0.87759637830000004),
    @("Spike-HA",   "PR AUC",    0.93083685640000002),
    @("RBD-NTD-S2", "Accuracy",  0.62198388576507502),
    @("RBD-NTD-S2", "Precision", 0.65714287757873502),
    @("RBD-NTD-S2", "Recall",    0.601307213306427),
    @("RBD-NTD-S2", "ROC AUC",   0.78563416004180897),
    @("RBD-NTD-S2", "PR AUC",    0.62366157770156805)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

$ws.Range("G16").Select()
